# Generate Report for Handoff
#
# The handback report tracks, per source file & target locale, the most
# recent handoff timestamp. The a9e883ab-dab4-40f8-b6f9-35275b0ef014.md
# file was just (re-)handed off, so its "Latest Handoff" timestamps are
# refreshed on the Overview sheet as well as on each per-locale detail
# sheet (row 5 in each case corresponds to that file).

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest Handoff Date" column (D) for the a9e883ab... row
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D5").Value = "2016-03-25 07:54:54"

# zh-cn detail sheet: "Latest Handoff Datetime" column (E) for the a9e883ab... row
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E5").Value = "2016-03-25 07:54:46"

# de-de detail sheet: "Latest Handoff Datetime" column (E) for the a9e883ab... row
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E5").Value = "2016-03-25 07:54:54"
